$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 11594
$ws.Range("E2").Value = 68
$ws.Range("F2").Value = 68
$ws.Range("G2").Value = 202
$ws.Range("H2").Value = 149
$ws.Range("I2").Value = 51
$ws.Range("J2").Value = 98
$ws.Range("K2").Value = 13797
$ws.Range("L2").Value = 2503
$ws.Range("M2").Value = 11294
$ws.Range("N2").Value = 6559
$ws.Range("O2").Value = 4734
$ws.Range("P2").Value = 185
$ws.Range("Q2").Value = 392
$ws.Range("R2").Value = -480
$ws.Range("S2").Value = -45
$ws.Range("T2").Value = 192
$ws.Range("U2").Value = 200
$ws.Range("V2").Value = 154
$ws.Range("W2").Value = 0.59
$ws.Range("X2").Value = 1.28
$ws.Range("Y2").Value = 0.77
$ws.Range("Z2").Value = 1.08
$ws.Range("AA2").Value = 22.17
$ws.Range("AB2").Value = 2519.03
$ws.Range("AC2").Value = 276
$ws.Range("AD2").Value = 35.7
$ws.Range("AE2").Value = 42884
$ws.Range("AF2").Value = 0.23
$ws.Range("AG2").Value = 180
$ws.Range("AH2").Value = 0.37
$ws.Range("AI2").Value = 54.07
$ws.Range("AJ2").Value = 18476380

# Row 3
$ws.Range("D3").Value = 10875
$ws.Range("E3").Value = 904
$ws.Range("F3").Value = 904
$ws.Range("G3").Value = 986
$ws.Range("H3").Value = 780
$ws.Range("I3").Value = 437
$ws.Range("J3").Value = 344
$ws.Range("K3").Value = 14071
$ws.Range("L3").Value = 2051
$ws.Range("M3").Value = 12020
$ws.Range("N3").Value = 7048
$ws.Range("O3").Value = 4973
$ws.Range("P3").Value = 185
$ws.Range("Q3").Value = 2016
$ws.Range("R3").Value = -1582
$ws.Range("S3").Value = -234
$ws.Range("T3").Value = 385
$ws.Range("U3").Value = 1631
$ws.Range("V3").Value = 13
$ws.Range("W3").Value = 8.31
$ws.Range("X3").Value = 7.18
$ws.Range("Y3").Value = 6.42
$ws.Range("Z3").Value = 5.6
$ws.Range("AA3").Value = 17.06
$ws.Range("AB3").Value = 2796.07
$ws.Range("AC3").Value = 2365
$ws.Range("AD3").Value = 5.26
$ws.Range("AE3").Value = 46076
$ws.Range("AF3").Value = 0.27
$ws.Range("AG3").Value = 180
$ws.Range("AH3").Value = 0.29
$ws.Range("AI3").Value = 6.3
$ws.Range("AJ3").Value = 18476380

# Row 4
$ws.Range("D4").Value = 10581
$ws.Range("E4").Value = 871
$ws.Range("F4").Value = 871
$ws.Range("G4").Value = 1006
$ws.Range("H4").Value = 779
$ws.Range("I4").Value = 459
$ws.Range("J4").Value = 320
$ws.Range("K4").Value = 14940
$ws.Range("L4").Value = 2200
$ws.Range("M4").Value = 12740
$ws.Range("N4").Value = 7494
$ws.Range("O4").Value = 5246
$ws.Range("P4").Value = 185
$ws.Range("Q4").Value = 1150
$ws.Range("R4").Value = -985
$ws.Range("S4").Value = -72
$ws.Range("T4").Value = 170
$ws.Range("U4").Value = 980
$ws.Range("V4").Value = 33
$ws.Range("W4").Value = 8.23
$ws.Range("X4").Value = 7.36
$ws.Range("Y4").Value = 6.31
$ws.Range("Z4").Value = 5.37
$ws.Range("AA4").Value = 17.27
$ws.Range("AB4").Value = 3034.76
$ws.Range("AC4").Value = 2483
$ws.Range("AD4").Value = 4.93
$ws.Range("AE4").Value = 48994
$ws.Range("AF4").Value = 0.25
$ws.Range("AG4").Value = 180
$ws.Range("AH4").Value = 0.29
$ws.Range("AI4").Value = 6
$ws.Range("AJ4").Value = 18476380

# Row 5
$ws.Range("D5").Value = 12213
$ws.Range("E5").Value = 974
$ws.Range("F5").Value = 974
$ws.Range("G5").Value = -89
$ws.Range("H5").Value = -72
$ws.Range("I5").Value = 129
$ws.Range("J5").Value = -201
$ws.Range("K5").Value = 14888
$ws.Range("L5").Value = 2300
$ws.Range("M5").Value = 12588
$ws.Range("N5").Value = 7601
$ws.Range("O5").Value = 4987
$ws.Range("P5").Value = 185
$ws.Range("Q5").Value = 1506
$ws.Range("R5").Value = -931
$ws.Range("S5").Value = -79
$ws.Range("T5").Value = 179
$ws.Range("U5").Value = 1327
$ws.Range("V5").Value = 47
$ws.Range("W5").Value = 7.97
$ws.Range("X5").Value = -0.59
$ws.Range("Y5").Value = 1.7
$ws.Range("Z5").Value = -0.49
$ws.Range("AA5").Value = 18.27
$ws.Range("AB5").Value = 3092.65
$ws.Range("AC5").Value = 697
$ws.Range("AD5").Value = 21.3
$ws.Range("AE5").Value = 49695
$ws.Range("AF5").Value = 0.3
$ws.Range("AG5").Value = 250
$ws.Range("AH5").Value = 0.34
$ws.Range("AI5").Value = 29.71
$ws.Range("AJ5").Value = 18476380

# Row 6
$ws.Range("D6").Value = 13148
$ws.Range("E6").Value = 624
$ws.Range("F6").Value = 624
$ws.Range("G6").Value = 580
$ws.Range("H6").Value = 363
$ws.Range("I6").Value = 254
$ws.Range("K6").Value = 15452
$ws.Range("L6").Value = 2661
$ws.Range("M6").Value = 12790
$ws.Range("N6").Value = 7789
$ws.Range("P6").Value = 185
$ws.Range("Q6").Value = 774
$ws.Range("R6").Value = -556
$ws.Range("S6").Value = -93
$ws.Range("T6").Value = 165
$ws.Range("U6").Value = 610
$ws.Range("V6").Value = 57
$ws.Range("W6").Value = 4.74
$ws.Range("X6").Value = 2.76
$ws.Range("Y6").Value = 3.29
$ws.Range("Z6").Value = 2.39
$ws.Range("AA6").Value = 20.8
$ws.Range("AB6").Value = 3200.05
$ws.Range("AC6").Value = 1372
$ws.Range("AD6").Value = 10.02
$ws.Range("AE6").Value = 50920
$ws.Range("AF6").Value = 0.27
$ws.Range("AG6").Value = 280
$ws.Range("AH6").Value = 2.04
$ws.Range("AI6").Value = 16.89
$ws.Range("AJ6").Value = 18476380

# Clear rows 7-9 data cells (D:AJ), keep A/B/C
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
